$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2017, 3468716000, 3720356000, 1100000000, 8333000),
    @(2018, 4469282000, 555913000, 2050000000, -3298000),
    @(2019, 5373593000, -1696013000, 2750000000, 114959000),
    @(2020, 6587000000, 2634000000, 3050000000, 73000000),
    @(2021, 8514000000, 1719000000, 3950000000, 96000000),
    @(2022, 9254000000, 868000000, 7068000000, 51000000),
    @(2023, 10291000000, 2833000000, 4400000000, -156000000)
)

$row = 2
foreach ($r in $data) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = [string]$r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
